$wb = $excel.ActiveWorkbook

# --- 1. Add the new "PAINCO" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PAINCO"

# --- 2. Header row (row 1), reusing the same header text as the other sheets ---
$ws.Range("A1").Value = "Peneiras`n (Mesh)"
$ws.Range("B1").Value = "wi `n (g)"
$ws.Range("C1").Value = "xi `n (%100)"
$ws.Range("D1").Value = "Xi `n (%100)"
$ws.Range("E1").Value = "- Di `n (μm)"
$ws.Range("F1").Value = "+ Di `n (μm)"
$ws.Range("G1").Value = "Di `n (μm)"
$ws.Range("H1").Value = "xi/Di"
$ws.Range("J1").Value = "Massa Inicial`n(g)"
$ws.Range("K1").Value = "dsp`n(μm)"

# --- 3. Sieve-size labels (column A) ---
$ws.Range("A2").Value = "-6+8"
$ws.Range("A3").Value = "-8+10"
$ws.Range("A4").Value = "-10+12"
$ws.Range("A5").Value = "-12+14"
$ws.Range("A6").Value = "-14+16"
$ws.Range("A7").Value = "-16+20"
$ws.Range("A8").Value = "-20+inf"

# --- 4. Mass data (column B) ---
$ws.Range("B2").Value = 0.4
$ws.Range("B3").Value = 532.8
$ws.Range("B4").Value = 61.8
$ws.Range("B5").Value = 2.5
$ws.Range("B6").Value = 1.0
$ws.Range("B7").Value = 0.5
$ws.Range("B8").Value = 0.3

# --- 5. Formula columns ---
$ws.Range("C2").Formula = '=B2/$J$2*100'
$ws.Range("C3").Formula = '=B3/$J$2*100'
$ws.Range("C4").Formula = '=B4/$J$2*100'
$ws.Range("C5").Formula = '=B5/$J$2*100'
$ws.Range("C6").Formula = '=B6/$J$2*100'
$ws.Range("C7").Formula = '=B7/$J$2*100'
$ws.Range("C8").Formula = '=B8/$J$2*100'

$ws.Range("D2").Formula = '=SUM(C3:C8)'
$ws.Range("D3").Formula = '=SUM(C4:C8)'
$ws.Range("D4").Formula = '=SUM(C4:C8)'
$ws.Range("D5").Formula = '=SUM(C5:C8)'
$ws.Range("D6").Formula = '=SUM(C7:C8)'
$ws.Range("D7").Formula = '=SUM(C8:C9)'
$ws.Range("D8").Formula = '=SUM(C9:C10)'

$ws.Range("E2").Value = 3360.0
$ws.Range("E3").Formula = '=F2'
$ws.Range("E4").Formula = '=F3'
$ws.Range("E5").Formula = '=F4'
$ws.Range("E6").Formula = '=F5'
$ws.Range("E7").Formula = '=F6'
$ws.Range("E8").Formula = '=F7'

$ws.Range("F2").Value = 2380.0
$ws.Range("F3").Value = 1680.0
$ws.Range("F4").Value = 1410.0
$ws.Range("F5").Value = 1190.0
$ws.Range("F6").Value = 1000.0
$ws.Range("F7").Value = 841.0
$ws.Range("F8").Value = 0.0

$ws.Range("G2").Formula = '=AVERAGE(E2:F2)'
$ws.Range("G3").Formula = '=AVERAGE(E3:F3)'
$ws.Range("G4").Formula = '=AVERAGE(E4:F4)'
$ws.Range("G5").Formula = '=AVERAGE(E5:F5)'
$ws.Range("G6").Formula = '=AVERAGE(E6:F6)'
$ws.Range("G7").Formula = '=AVERAGE(E7:F7)'
$ws.Range("G8").Formula = '=AVERAGE(E8:F8)'

$ws.Range("H2").Formula = '=C2/100/G2'
$ws.Range("H3").Formula = '=C3/100/G3'
$ws.Range("H4").Formula = '=C4/100/G4'
$ws.Range("H5").Formula = '=C5/100/G5'
$ws.Range("H6").Formula = '=C6/100/G6'
$ws.Range("H7").Formula = '=C7/100/G7'
$ws.Range("H8").Formula = '=C8/100/G8'

$ws.Range("J2").Formula = '=SUM(B2:B8)'
$ws.Range("K2").Formula = '=SUM(C2:C8)/100/SUM(H2:H8)'

Write-Host ("Done populating PAINCO sheet")
